$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add the two new rows of the Key -> Mapping table
$ws.Range("A28").Value = "MQ 6.1"
$ws.Range("B28").Value = "Total"

$ws.Range("A29").Value = "MQ 5.1"
$ws.Range("B29").Value = "Total"

# Update selection / view state to match the saved workbook
$ws.Range("A29").Select()
